$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (a "type" row) - shifts existing rows 2..5 down to 3..6
$ws.Rows.Item(2).Insert()

# Fill in the new type row (row 2) describing the data type of each column
$ws.Range("A2").Value = "int"
$ws.Range("B2").Value = "string"
$ws.Range("C2").Value = "string"
$ws.Range("D2").Value = "int"
$ws.Range("E2").Value = "int"
$ws.Range("F2").Value = "int"
$ws.Range("G2").Value = "int"
$ws.Range("H2").Value = "int"
$ws.Range("I2").Value = "int"
$ws.Range("J2").Value = "int"
$ws.Range("K2").Value = "int[]"
$ws.Range("L2").Value = "int"
$ws.Range("M2").Value = "int[]"
$ws.Range("N2").Value = "float"
$ws.Range("O2").Value = "float"

# Update the first data row (now row 4, previously row 3): SkillIds becomes a bracketed
# list, and the breathing rate becomes a text value with an 'f' suffix
$ws.Range("K4").Value = "[1,2,3]"
$ws.Range("O4").Value = "4.8f"

# Update the selection to reflect the commit's final cursor position
$ws.Range("O4").Select()
$excel.ActiveWindow.ScrollColumn = 3
